{"js": "// Bump the version string in the title block from \"V3.3.6.2\" to \"V3.3.7.2\".\n// The paragraph originally holds the version number as four separate runs\n// (\"V3.3.\", \"6\", \".\", \"2\") plus a trailing _GoBack bookmark; only the \"6\"\n// run's text actually needs to change to \"7\", but re-typing that digit also\n// shifts Word's _GoBack bookmark (which tracks the most recent edit point)\n// to the front of the paragraph, so we reproduce that too.\n\nconst body = context.document.body;\n\n// Locate the version-number paragraph via the stable \"V3.3.\" prefix.\nconst hits = body.search(\"V3.3.\", { matchCase: true });\nhits.load(\"items\");\nawait context.sync();\n\nif (hits.items.length === 0) {\n  throw new Error(\"Could not find the 'V3.3.' version string\");\n}\n\nconst versionHit = hits.items[0];\nconst para = versionHit.paragraphs.getFirst();\nconst paraRange = para.getRange();\n\n// Find the lone \"6\" that sits between \"V3.3.\" and \".2\" inside that paragraph.\nconst digitHits = paraRange.search(\"6\", { matchCase: true });\ndigitHits.load(\"items\");\nawait context.sync();\n\nif (digitHits.items.length !== 1) {\n  throw new Error(\"Expected exactly one '6' in the version paragraph, found \" + digitHits.items.length);\n}\n\n// Replace just that run's text; because the match exactly spans the run,\n// the surrounding runs (\"V3.3.\" and \".\") are left untouched.\ndigitHits.items[0].insertText(\"7\", \"Replace\");\nawait context.sync();\n\n// Re-anchor the _GoBack bookmark at the start of the paragraph, matching\n// Word's behaviour of moving it to the location of the latest edit.\ncontext.document.deleteBookmark(\"_GoBack\");\nconst startRange = para.getRange(\"Start\");\nstartRange.insertBookmark(\"_GoBack\");\n\nawait context.sync();\n", "ps1": "# Bump the version string in the title block from \"V3.3.6.2\" to \"V3.3.7.2\".\n# The paragraph originally holds the version number as four separate runs\n# (\"V3.3.\", \"6\", \".\", \"2\") plus a trailing _GoBack bookmark; only the \"6\"\n# run's text actually needs to change to \"7\", but re-typing that digit also\n# shifts Word's _GoBack bookmark (which tracks the most recent edit point)\n# to the front of the paragraph, so we reproduce that too.\n\n$d = $word.ActiveDocument\n\n# Locate the version-number paragraph via the stable \"V3.3.\" prefix, then\n# expand the find hit out to the whole paragraph.\n$paraRange = $d.Content\n$paraRange.Find.Execute(\"V3.3.\")\n$paraRange.Expand(4)\n\n# Find the lone \"6\" within that paragraph (work on a duplicate so paraRange\n# itself keeps spanning the whole paragraph).\n$digitRange = $paraRange.Duplicate\n$digitRange.Find.Execute(\"6\")\n\n# Toggle Bold on/off around the text replace so the edited run does not\n# silently re-merge with its identically-formatted neighbouring run; the\n# final Bold=0 restores the original formatting, leaving the run separate.\n$digitRange.Bold = 1\n$digitRange.Text = \"7\"\n$digitRange.Bold = 0\n\n# Re-anchor the _GoBack bookmark at the start of the paragraph, matching\n# Word's behaviour of moving it to the location of the latest edit.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n  $d.Bookmarks(\"_GoBack\").Delete()\n}\n$startRange = $paraRange.Duplicate\n$startRange.Collapse(1)\n$d.Bookmarks.Add(\"_GoBack\", $startRange)\n"}
